$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "PHOR"
$ws.Range("B8").Value = "ФосАгро"
$ws.Range("C8").Value = "9978b56f-782a-4a80-a4b1-a48cbecfd194"
$ws.Range("D8").Value = "PHOR"
$ws.Range("E8").Value = "RU000A0JRKT8"
$ws.Range("F8").Value = "BBG004S689R0"
$ws.Range("G8").Value = "TQBR"
$ws.Range("H8").Value = "share"
$ws.Range("I8").Value = "5a3d1efd-f8a0-478e-a10e-bb7f990f9c87"

$ws.Range("A9").Value = "HEAD"
$ws.Range("B9").Value = "Хэдхантер"
$ws.Range("C9").Value = "3fe80143-1313-42eb-9884-5d68b39e265e"
$ws.Range("D9").Value = "HEAD"
$ws.Range("E9").Value = "RU000A107662"
$ws.Range("F9").Value = "TCS20A107662"
$ws.Range("G9").Value = "TQBR"
$ws.Range("H9").Value = "share"
$ws.Range("I9").Value = "911552ef-a892-4b33-9df1-c0d6c4a2307d"
